$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before column N (14th column), shifting Late/heading/Outstanding right.
$ws.Range("N1").EntireColumn.Insert()

# Give the newly inserted column N a custom (non-autofit) width, close to the
# "In Advance"/"Late" column width it sits next to. The previously existing
# columns (now shifted to O, P, Q) keep their original widths automatically
# because Insert() shifts the <col> definitions along with the data.
$ws.Range("N1").EntireColumn.ColumnWidth = 9.8

$ws.Range("R6").Select()
